# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (a copy of the "2021-Q4" sheet, which
#    already has the right fund-holding table layout/styles) positioned
#    right before the "总计" (Total) summary sheet, then overwrite its
#    contents with the new 2022-Q1 fund holder data.
# 2) Insert a new top data row into the "总计" sheet for the 2022-Q1 summary
#    and bump the existing running index column down by one.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# --- Step 1: create the "2022-Q1" sheet just before "总计" ---
$q4Sheet.Copy($totalSheet)
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($totalSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# The copied "2021-Q4" sheet has 13 data rows (rows 2-14); the new data only
# has 10 data rows (rows 2-11), so drop the trailing 3 rows.
$newSheet.Rows.Item(12).Delete()
$newSheet.Rows.Item(12).Delete()
$newSheet.Rows.Item(12).Delete()

# Header row (text already, just overwrite the labels in place)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @(0, "003567", "华夏行业景气混合", "112.21", "91.63", "2.80", "3.1419", 6),
    @(1, "162212", "泰达宏利红利先锋混合", "0.68", "92.79", "3.96", "0.0269", 7),
    @(2, "010703", "财通智选消费股票A", "1.12", "92.43", "2.34", "0.0262", 9),
    @(3, "002145", "诺安景鑫灵活配置混合", "0.53", "83.45", "4.59", "0.0243", 2),
    @(4, "010704", "财通智选消费股票C", "0.54", "92.43", "2.34", "0.0126", 9),
    @(5, "006818", "安信盈利驱动股票A", "0.27", "83.93", "4.62", "0.0125", 7),
    @(6, "006819", "安信盈利驱动股票C", "0.17", "83.93", "4.62", "0.0079", 7),
    @(7, "001657", "长安鑫富领先灵活配置混合", "0.07", "30.32", "1.97", "0.0014", 9),
    @(8, "005537", "中航新起航灵活配置混合A", "0.03", "87.09", "4.19", "0.0013", 9),
    @(9, "005538", "中航新起航灵活配置混合C", "0.01", "87.09", "4.19", "0.0004", 9)
)

# Text-valued columns (code, scale, total position, position ratio, market
# value) must stay text (e.g. "003567" keeps its leading zeros) rather than
# being auto-coerced into numbers. Force text format while writing, then
# clear the formatting override again so no stray style survives.
$textCols = @("B", "D", "E", "F", "G")

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

foreach ($col in $textCols) {
    $range = $newSheet.Range("$col" + "2:" + "$col" + "11")
    $range.NumberFormat = "@"
}

foreach ($row in $fundRows) {
    $rowIndex = $row[0] + 2
    $newSheet.Range("B$rowIndex").Value = $row[1]
    $newSheet.Range("D$rowIndex").Value = $row[3]
    $newSheet.Range("E$rowIndex").Value = $row[4]
    $newSheet.Range("F$rowIndex").Value = $row[5]
    $newSheet.Range("G$rowIndex").Value = $row[6]
}

foreach ($col in $textCols) {
    $range = $newSheet.Range("$col" + "2:" + "$col" + "11")
    $range.ClearFormats()
}

# --- Step 2: add the 2022-Q1 summary row atop the "总计" sheet ---
$totalSheet.Rows.Item(2).Insert()

# Reuse row 3's existing cell formatting for the new row 2 (copy+paste the
# format only, so the new A2 correctly picks up the same style used by the
# rest of column A instead of a brand new style record).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 10
$totalSheet.Range("D2").Value = 3.26

# Bump the running index in column A for the rows that used to start at 0
# (now rows 3-7) so the sequence reads 0,1,2,3,4,5 top to bottom.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
